$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data set is a rolling window: the oldest reading (row 15) is dropped,
# rows 16-18 shift up to rows 15-17, and a brand-new reading is appended
# at row 17. Net effect: row 18 disappears, dimension shrinks to A1:J17.

# Row 15 <- old row 16 values
$ws.Range("A15").Value = 45874.54183811343
$ws.Range("D15").Value = 19.78
$ws.Range("E15").Value = 76.77
$ws.Range("F15").Value = 578.3200000000001
$ws.Range("G15").Value = 14.38
$ws.Range("J15").Value = "13:00:14"

# Row 16 <- old row 17 values
$ws.Range("A16").Value = 45874.5836159375
$ws.Range("D16").Value = 20.56
$ws.Range("E16").Value = 74.25
$ws.Range("F16").Value = 82.62
$ws.Range("G16").Value = 13.87
$ws.Range("J16").Value = "14:00:24"

# Row 17 <- brand new reading
$ws.Range("A17").Value = 45874.66686496891
$ws.Range("D17").Value = 21.5
$ws.Range("E17").Value = 71.73
$ws.Range("F17").Value = 298.83
$ws.Range("G17").Value = 11.38
$ws.Range("J17").Value = "16:00:17"

# Row 18 no longer exists in the data set - delete it entirely
$ws.Rows("18").Delete()
